$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (all stored as text/inline strings)
$changes = @{
    "D2"  = "331.83"
    "E2"  = "1.94%"
    "D3"  = "41.13"
    "E3"  = "3.19%"
    "D4"  = "5.749"
    "E4"  = "-2.03%"
    "D5"  = "0.08160"
    "E5"  = "2.02%"
    "E6"  = "7.58%"
    "D7"  = "8.748"
    "E7"  = "0.42%"
    "E8"  = "-1.55%"
    "E9"  = "-0.19%"
    "D10" = "0.9207"
    "D11" = "0.1244"
    "E11" = "-1.22%"
    "D12" = "0.1950"
    "E12" = "-0.69%"
    "D13" = "8.303"
    "E13" = "-6.07%"
    "D14" = "0.09377"
    "E14" = "2.12%"
    "D15" = "0.03654"
    "E15" = "1.96%"
    "E16" = "9.52%"
    "D17" = "0.001295"
    "E17" = "-0.94%"
    "D18" = "0.006170"
    "E18" = "0.99%"
    "D19" = "3.385"
    "E19" = "0.91%"
    "E20" = "-1.19%"
    "E21" = "-1.13%"
    "E22" = "9.53%"
    "D23" = "0.04426"
    "E23" = "-0.39%"
    "D24" = "0.001260"
    "E24" = "-0.10%"
    "D25" = "0.004341"
    "E25" = "-0.01%"
    "E26" = "8.40%"
    "D39" = "0.02776"
    "E39" = "14.42%"
    "D40" = "0.05511"
    "E40" = "4.79%"
    "D41" = "0.007638"
    "E41" = "1.74%"
    "D42" = "0.009941"
    "E42" = "14.35%"
    "E43" = "0.85%"
    "D44" = "0.002119"
    "E44" = "0.64%"
    "E45" = "12.41%"
    "D46" = "0.00006724"
    "E46" = "-1.82%"
    "E47" = "-0.45%"
    "E48" = "59.74%"
    "E49" = "1.23%"
    "E50" = "-0.45%"
    "E51" = "-0.45%"
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$addr]
}
